# Generate Report for handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the ae38855b-...-7ca250d3... entry (row 2) on both the
# "zh-cn" and "de-de" worksheets to reflect the new report generation time.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-20 07:50:07"
$wsZhCn.Range("G2").Value = "2016-01-20 07:50:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-20 07:50:19"
$wsDeDe.Range("G2").Value = "2016-01-20 07:51:13"
